$wb = $excel.ActiveWorkbook

# --- Sheet "Data": append a new observation row (row 96)
#     date 2023-07-26 (serial 45133), value 8243.343999999999
$dataSheet = $wb.Worksheets.Item("Data")

# Duplicate the date-cell formatting from the row above so the new date
# cell gets the same style as the rest of the date column.
$dataSheet.Range("A95").Copy() | Out-Null
$dataSheet.Range("A96").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$dataSheet.Cells.Item(96, 1).Value = 45133
$dataSheet.Cells.Item(96, 2).Value = 8243.343999999999

# --- Sheet "SeriesInfo": refresh metadata to match the new pull
$infoSheet = $wb.Worksheets.Item("SeriesInfo")

function Set-TextValue($range, [string]$text) {
    # Force the written value to remain text instead of letting Excel's
    # autodetection turn date-looking strings into date serials.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $infoSheet.Range("B3")  "2023-08-03"
Set-TextValue $infoSheet.Range("B4")  "2023-08-03"
Set-TextValue $infoSheet.Range("B7")  "2023-07-26"
Set-TextValue $infoSheet.Range("B14") "2023-07-27 15:33:02-05"

$infoSheet.Range("B15").Value = 95
